# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1: header cell "Save" - copy formatting from the neighboring header cell (G1)
# so it picks up the same bold/border/centered style used by the other headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# H2: data cell with the save value for this row
$ws.Range("H2").Value = 1
